$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K21").Value = "revision milestone"
$ws.Range("J9").Value = "estimated time"
$ws.Range("K9").Value = "actual time"
$ws.Range("L9").Value = "Milestones"
$ws.Range("M21").Value = "project libre"

Write-Host "done"
